# Adds a new "Q9" (K-column) data point across the ifoCAST-matched error series.
# Rows 2-14 gain a 10th quarter value; rows 15-24 lose their last (oldest) quarter
# value as the rolling window advances by one column.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","G","H","I","J","K")

$rowData = @{}
$rowData[2] = @(2.489322709745834, 9.958473040745833, -8.212788609254167, -0.07311760925416594, 2.074453258745834, -1.365576609254166, -1.426750609254166, 0.6147053907458341, 0.03211039074583408, 0.1277855347458341)
$rowData[3] = @(9.643547872324863, -8.527713777675137, -0.3880427776751375, 1.759528090324862, -1.680501777675137, -1.741675777675137, 0.2997802223248625, -0.2828147776751375, -0.1871396336751375, -0.4475587776751375)
$rowData[4] = @(-9.584088737857668, -1.444417737857667, 0.7031531301423328, -2.736876737857667, -2.798050737857667, -0.7565947378576671, -1.339189737857667, -1.243514593857667, -1.503933737857667, -0.6684787378576671)
$rowData[5] = @(-0.5737713618867365, 1.573799506113263, -1.866230361886736, -1.927404361886736, 0.1140516381132635, -0.4685433618867365, -0.3728682178867365, -0.6332873618867365, 0.2021676381132635, 0.1020946851132635)
$rowData[6] = @(1.56314828978338, -1.87688157821662, -1.93805557821662, 0.10340042178338, -0.47919457821662, -0.38351943421662, -0.64393857821662, 0.19151642178338, 0.09144346878338, -0.26234957821662)
$rowData[7] = @(-1.985496110432761, -2.046670110432761, -0.005214110432760821, -0.5878091104327609, -0.4921339664327609, -0.7525531104327607, 0.0829018895672392, -0.01717106343276081, -0.3709641104327608, -0.2227861104327608)
$rowData[8] = @(-1.838568389465022, 0.2028876105349782, -0.3797073894650218, -0.2840322454650218, -0.5444513894650218, 0.2910036105349783, 0.1909306575349782, -0.1628623894650218, -0.01468438946502179, 0.1427616105349782)
$rowData[9] = @(0.3193857809087786, -0.2632092190912214, -0.1675340750912214, -0.4279532190912214, 0.4075017809087786, 0.3074288279087786, -0.04636421909122137, 0.1018137809087786, 0.2592597809087786, -0.2414032190912214)
$rowData[10] = @(-0.2177158050940863, -0.1220406610940863, -0.3824598050940862, 0.4529951949059138, 0.3529222419059138, -0.0008708050940862688, 0.1473071949059137, 0.3047531949059137, -0.1959098050940863, 0.08445119490591374)
$rowData[11] = @(-0.1395947824540507, -0.4000139264540507, 0.4354410735459493, 0.3353681205459493, -0.01842492645405072, 0.1297530735459493, 0.2871990735459493, -0.2134639264540507, 0.06689707354594927, -0.1050779264540507)
$rowData[12] = @(-0.3119066075030444, 0.5235483924969555, 0.4234754394969555, 0.06968239249695551, 0.2178603924969555, 0.3753063924969555, -0.1253566075030445, 0.1550043924969555, -0.01697060750304449, 0.2888923924969555)
$rowData[13] = @(0.7021230259847391, 0.6020500729847391, 0.2482570259847391, 0.3964350259847391, 0.5538810259847391, 0.0532180259847391, 0.3335790259847391, 0.1616040259847391, 0.4674670259847391, -0.1448629740152609)
$rowData[14] = @(1.514077300737389, 1.160284253737389, 1.308462253737389, 1.465908253737389, 0.9652452537373891, 1.245606253737389, 1.073631253737389, 1.379494253737389, 0.7671642537373891, 1.455535253737389)
$rowData[15] = @(0.2163100177716323, 0.3644880177716323, 0.5219340177716323, 0.02127101777163229, 0.3016320177716323, 0.1296570177716323, 0.4355200177716323, -0.1768099822283677, 0.5115610177716323, $null)
$rowData[16] = @(0.3684559152847414, 0.5259019152847414, 0.0252389152847414, 0.3055999152847414, 0.1336249152847414, 0.4394879152847414, -0.1728420847152586, 0.5155289152847414, $null, $null)
$rowData[17] = @(0.6615420054549828, 0.1608790054549828, 0.4412400054549828, 0.2692650054549828, 0.5751280054549828, -0.03720199454501721, 0.6511690054549828, $null, $null, $null)
$rowData[18] = @(-0.07992717640068192, 0.2004338235993181, 0.02845882359931809, 0.3343218235993181, -0.2780081764006819, 0.4103628235993181, $null, $null, $null, $null)
$rowData[19] = @(0.1551026083472217, -0.01687239165277828, 0.2889906083472217, -0.3233393916527783, 0.3650316083472218, $null, $null, $null, $null, $null)
$rowData[20] = @(-0.08373381358883215, 0.2221291864111679, -0.3902008135888321, 0.2981701864111679, $null, $null, $null, $null, $null, $null)
$rowData[21] = @(0.1925429840544512, -0.4197870159455487, 0.2685839840544513, $null, $null, $null, $null, $null, $null, $null)
$rowData[22] = @(-0.4379376726751896, 0.2504333273248104, $null, $null, $null, $null, $null, $null, $null, $null)
$rowData[23] = @(0.2324015030761111, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$rowData[24] = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $cellRef = "$($cols[$i])$r"
        if ($null -eq $vals[$i]) {
            $ws.Range($cellRef).ClearContents()
        } else {
            $ws.Range($cellRef).Value = $vals[$i]
        }
    }
}

Write-Output "Updated $($rowData.Keys.Count) rows with ifoCAST full series evaluation data"
